# "Generate Report for Handback"
#
# The handback run for both localized files (c900173d-... and
# e284f25b-...) has completed for the zh-cn and de-de languages:
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   - The "Latest Target File" (F) / "Latest Handback File" (G) columns
#     get populated (same file names/links as the source .md and the
#     handoff .xlf, since the handback round-tripped the exact files).
#   - The "Latest Handback DateTime" (H) moves from the zero-date sentinel
#     to the real handback timestamp.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"
$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/3f6362359bf40a7c2ec42e80284291a9442bc558/e2e/"

# ---- Overview sheet: refresh the mirrored status cells ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# ---- zh-cn sheet ----
$ws = $wb.Worksheets.Item("zh-cn")
$xlfBaseUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/74f2b94c4e07e6852d41d28ce0fa5e1476689e6c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/"
$handbackDt = "2016-03-12 04:47:30"

$mdName = "c900173d-dd2c-40e4-9e4f-0c102b4800c9.md"
$xlfName = "c900173d-dd2c-40e4-9e4f-0c102b4800c9.5b69022692ffc2765573cf5bf61fb24df16cea19.zh-cn.xlf"
$ws.Range("C2").Value = $newStatus
$ws.Range("F2").Value = $mdName
$ws.Hyperlinks.Add($ws.Range("F2"), ($mdUrl + $mdName), "", "", $mdName)
$ws.Range("G2").Value = $xlfName
$ws.Hyperlinks.Add($ws.Range("G2"), ($xlfBaseUrl + $xlfName), "", "", $xlfName)
$ws.Range("H2").Value = $handbackDt

$mdName = "e284f25b-4ebb-441f-b441-a5d8ef7b49a3.md"
$xlfName = "e284f25b-4ebb-441f-b441-a5d8ef7b49a3.ccb4d763564db2d9fcd1bbab751c0edfe862d7bd.zh-cn.xlf"
$ws.Range("C3").Value = $newStatus
$ws.Range("F3").Value = $mdName
$ws.Hyperlinks.Add($ws.Range("F3"), ($mdUrl + $mdName), "", "", $mdName)
$ws.Range("G3").Value = $xlfName
$ws.Hyperlinks.Add($ws.Range("G3"), ($xlfBaseUrl + $xlfName), "", "", $xlfName)
$ws.Range("H3").Value = $handbackDt

# ---- de-de sheet ----
$ws = $wb.Worksheets.Item("de-de")
$xlfBaseUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/70de055e5b088e026e59ecd8430585d491450fb9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/"
$handbackDt = "2016-03-12 04:47:36"

$mdName = "c900173d-dd2c-40e4-9e4f-0c102b4800c9.md"
$xlfName = "c900173d-dd2c-40e4-9e4f-0c102b4800c9.5b69022692ffc2765573cf5bf61fb24df16cea19.de-de.xlf"
$ws.Range("C2").Value = $newStatus
$ws.Range("F2").Value = $mdName
$ws.Hyperlinks.Add($ws.Range("F2"), ($mdUrl + $mdName), "", "", $mdName)
$ws.Range("G2").Value = $xlfName
$ws.Hyperlinks.Add($ws.Range("G2"), ($xlfBaseUrl + $xlfName), "", "", $xlfName)
$ws.Range("H2").Value = $handbackDt

$mdName = "e284f25b-4ebb-441f-b441-a5d8ef7b49a3.md"
$xlfName = "e284f25b-4ebb-441f-b441-a5d8ef7b49a3.ccb4d763564db2d9fcd1bbab751c0edfe862d7bd.de-de.xlf"
$ws.Range("C3").Value = $newStatus
$ws.Range("F3").Value = $mdName
$ws.Hyperlinks.Add($ws.Range("F3"), ($mdUrl + $mdName), "", "", $mdName)
$ws.Range("G3").Value = $xlfName
$ws.Hyperlinks.Add($ws.Range("G3"), ($xlfBaseUrl + $xlfName), "", "", $xlfName)
$ws.Range("H3").Value = $handbackDt

Write-Output "Handback report generated."
